$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List of artifacts")

for ($row = 2; $row -le 6; $row++) {
    $cell = $ws.Cells.Item($row, 5)
    $oldFormula = $cell.Formula
    $newFormula = $oldFormula.Replace('"$TRG_TBL",D' + $row, '"$TRG_TBL",A' + $row)
    $cell.Formula = $newFormula
}

$ws.Range("E6").Select()
